$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (contain only digits/dot) need the
# cell pre-formatted as Text so Excel keeps them as strings (preserving
# formatting such as trailing zeros) instead of auto-converting to numbers.
$textValues = @{
    "D5" = '266.08'
    "D6" = '0.622'
    "D8" = '57.35'
    "D9" = '0.388'
    "D10" = '0.0787'
    "D12" = '14.64'
    "D14" = '0.819'
    "D15" = '21.07'
    "D16" = '5.33'
    "D19" = '70.31'
    "D21" = '5.24'
    "D22" = '229.94'
    "D23" = '2.73'
    "D25" = '2.34'
    "D26" = '164.77'
    "D27" = '9.09'
    "D28" = '19.84'
    "D29" = '0.129'
    "D30" = '1.37'
    "D32" = '0.0665'
    "D33" = '4.68'
    "D34" = '4.57'
    "D35" = '2.45'
    "D36" = '1.81'
    "D37" = '3.39'
    "D39" = '5.30'
    "D40" = '3.07'
    "D41" = '1.24'
    "D42" = '0.0953'
    "D43" = '0.0216'
    "D45" = '91.94'
    "D46" = '16.02'
    "D47" = '1.05'
    "D48" = '7.21'
    "D49" = '2.89'
    "D50" = '2.00'
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $textValues[$addr]
    $ws.Range($addr).NumberFormat = "General"
}

# Remaining cells (percentages, links, coin names) can be set directly.
$plainValues = @{
    "D2" = '37.590.55'
    "E2" = '  +0.66%  '
    "D3" = '2.032.09'
    "E3" = '  +1.35%  '
    "E4" = '  -0.12%  '
    "E5" = '  +7.70%  '
    "E6" = '  -1.01%  '
    "E7" = '  -0.09%  '
    "E8" = '  -5.65%  '
    "E9" = '  +1.47%  '
    "E10" = '  -1.80%  '
    "E11" = '  -1.57%  '
    "E12" = '  -2.67%  '
    "D13" = '2.330.32'
    "E13" = '  +1.33%  '
    "E14" = '  -3.61%  '
    "E15" = '  -8.23%  '
    "E16" = '  -2.38%  '
    "D17" = '2.042.33'
    "E17" = '  +1.40%  '
    "D18" = '37.523.56'
    "E18" = '  +0.64%  '
    "E19" = '  -0.38%  '
    "D20" = '0.0₃0851'
    "E20" = '  -1.94%  '
    "E21" = '  +0.57%  '
    "E22" = '  -0.54%  '
    "E23" = '  +7.54%  '
    "E24" = '  -0.08%  '
    "E25" = '  -1.14%  '
    "E26" = '  +0.68%  '
    "E27" = '  -3.43%  '
    "E28" = '  +0.73%  '
    "E29" = '  -10.11%  '
    "E30" = '  +2.27%  '
    "E31" = '  -0.72%  '
    "E32" = '  +1.93%  '
    "E33" = '  -3.77%  '
    "E34" = '  +0.41%  '
    "E35" = '  +2.74%  '
    "B36" = 'WEMIXToken'
    "C36" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "E36" = '  +0.63%  '
    "B37" = 'RenderToken'
    "C37" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    "E37" = '  +1.93%  '
    "E38" = '  -0.22%  '
    "E39" = '  -4.43%  '
    "E40" = '  +4.94%  '
    "E41" = '  +3.81%  '
    "E42" = '  -2.84%  '
    "E43" = '  +0.65%  '
    "D44" = '1.421.37'
    "E44" = '  +3.16%  '
    "E45" = '  +1.00%  '
    "E46" = '  -3.59%  '
    "E47" = '  +0.08%  '
    "E48" = '  -1.06%  '
    "E49" = '  +1.29%  '
    "E50" = '  -1.54%  '
    "D51" = '2.221.84'
    "E51" = '  +1.32%  '
}

foreach ($addr in $plainValues.Keys) {
    $ws.Range($addr).Value = $plainValues[$addr]
}
